$d = $word.ActiveDocument

$pairs = @(
    @{old="87÷3="; new="29÷4="},
    @{old="18÷7="; new="25÷5="},
    @{old="20÷2="; new="50÷7="},
    @{old="59÷5="; new="33÷6="},
    @{old="63÷7="; new="17÷8="},
    @{old="58÷3="; new="81÷3="},
    @{old="63÷9="; new="45÷8="},
    @{old="57÷6="; new="48÷7="},
    @{old="43÷5="; new="48÷5="},
    @{old="45÷2="; new="43÷8="},
    @{old="60÷6="; new="79÷5="},
    @{old="11÷9="; new="30÷6="},
    @{old="90÷7="; new="49÷6="},
    @{old="62÷5="; new="63÷8="},
    @{old="56÷3="; new="86÷6="},
    @{old="30÷5="; new="93÷6="},
    @{old="88÷3="; new="60÷8="},
    @{old="88÷4="; new="48÷6="},
    @{old="13÷2="; new="16÷8="},
    @{old="26÷8="; new="41÷4="},
    @{old="52÷4="; new="49÷8="},
    @{old="57÷7="; new="51÷5="},
    @{old="74÷7="; new="59÷4="},
    @{old="39÷7="; new="27÷6="},
    @{old="86÷9="; new="76÷6="}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
